# Update OCI_Based_Unrealized_Losses_to_Assets (J) and
# OCI_Based_Unrealized_Losses_to_Equity (K) values for the bank indicators
# table, as part of adding the Spearman correlation matrix / stress-test
# rules functionality.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Value = 0.03
$ws.Range("K2").Value = 1.31

$ws.Range("J3").Value = 0.07000000000000001
$ws.Range("K3").Value = 2.51

$ws.Range("K4").Value = 5.56

$ws.Range("J5").Value = 0.06
$ws.Range("K5").Value = 1.25

$ws.Range("J6").Value = 0.03
$ws.Range("K6").Value = 1.17

$ws.Range("J8").Value = 0.08
$ws.Range("K8").Value = 2.34

$ws.Range("J9").Value = 0.09
$ws.Range("K9").Value = 3.46

$ws.Range("J10").Value = 0.01
$ws.Range("K10").Value = 0.38
